$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G20").Value = 2.35
$ws.Range("H20").Value = 3.1
$ws.Range("I20").Value = 2.92
$ws.Range("M20").Value = 2.4
$ws.Range("AE20").Value = 7.3
$ws.Range("AH20").Value = 35
$ws.Range("H21").Value = 3.6
$ws.Range("L21").Value = 1.47
$ws.Range("M21").Value = 2.35
$ws.Range("R21").Value = 2.45
$ws.Range("S21").Value = 1.42
$ws.Range("X21").Value = 16.5
$ws.Range("Z21").Value = 6.8
$ws.Range("Z26").Value = 10
$ws.Range("AB26").Value = 15
$ws.Range("AF26").Value = 23
$ws.Range("G32").Value = 2.37
$ws.Range("I32").Value = 4.1
$ws.Range("K32").Value = 4.1
$ws.Range("O32").Value = 1.32
$ws.Range("P32").Value = 1.7
$ws.Range("Q32").Value = 2.02
$ws.Range("T32").Value = 5
$ws.Range("U32").Value = 9.75
$ws.Range("V32").Value = 10.25
$ws.Range("W32").Value = 27
$ws.Range("X32").Value = 28
$ws.Range("Z32").Value = 4.1
$ws.Range("AE32").Value = 7.4
$ws.Range("AF32").Value = 21
$ws.Range("AG32").Value = 15
$ws.Range("AH32").Value = 80
$ws.Range("AI32").Value = 60
$ws.Range("AJ32").Value = 80
$ws.Range("G38").Value = 2.7
$ws.Range("I38").Value = 2.4
$ws.Range("N38").Value = 2
$ws.Range("P38").Value = 1.38
$ws.Range("Q38").Value = 2.47
$ws.Range("T38").Value = 6.8
$ws.Range("U38").Value = 10.75
$ws.Range("V38").Value = 8.5
$ws.Range("W38").Value = 24
$ws.Range("X38").Value = 19.5
$ws.Range("Y38").Value = 27
$ws.Range("Z38").Value = 8.25
$ws.Range("AE38").Value = 6.4
$ws.Range("AF38").Value = 9.5
$ws.Range("AG38").Value = 8
$ws.Range("AH38").Value = 20
$ws.Range("AI38").Value = 16.5
$ws.Range("AJ38").Value = 25
$ws.Range("G40").Value = 2.7
$ws.Range("I40").Value = 2.63
$ws.Range("T40").Value = 7
$ws.Range("U40").Value = 12
$ws.Range("W40").Value = 29
$ws.Range("AE40").Value = 6.5
$ws.Range("AH40").Value = 26
$ws.Range("G42").Value = 1.85
$ws.Range("H42").Value = 4.2
$ws.Range("I42").Value = 3.3
$ws.Range("J42").Value = 1.02
$ws.Range("K42").Value = 10
$ws.Range("L42").Value = 1.12
$ws.Range("M42").Value = 5.3
$ws.Range("N42").Value = 1.4
$ws.Range("O42").Value = 2.75
$ws.Range("P42").Value = 1.23
$ws.Range("Q42").Value = 3.8
$ws.Range("R42").Value = 1.42
$ws.Range("S42").Value = 2.65
$ws.Range("T42").Value = 13
$ws.Range("W42").Value = 18
$ws.Range("X42").Value = 12.5
$ws.Range("Y42").Value = 17
$ws.Range("Z42").Value = 10
$ws.Range("AA42").Value = 9.25
$ws.Range("AE42").Value = 18
$ws.Range("AF42").Value = 24
$ws.Range("AG42").Value = 12
$ws.Range("AH42").Value = 45
$ws.Range("AI42").Value = 24
$ws.Range("AJ42").Value = 23
$ws.Range("R46").Value = 1.57
$ws.Range("S46").Value = 2.12
$ws.Range("T46").Value = 8.75
$ws.Range("U46").Value = 11.75
$ws.Range("Y46").Value = 24
$ws.Range("AA46").Value = 6.5
$ws.Range("AB46").Value = 12
$ws.Range("AE46").Value = 11.25
$ws.Range("AF46").Value = 18
$ws.Range("AJ46").Value = 27
$ws.Range("G48").Value = 2.5
$ws.Range("I48").Value = 2.6
$ws.Range("P48").Value = 1.25
$ws.Range("Q48").Value = 3.75
$ws.Range("W48").Value = 26
$ws.Range("AG48").Value = 10
$ws.Range("AH48").Value = 26
$ws.Range("K50").Value = 19
$ws.Range("L50").Value = 1.13
$ws.Range("M50").Value = 6
$ws.Range("N50").Value = 1.44
$ws.Range("O50").Value = 2.7
